$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "FDR corrected p-value" column (old K); "result" (old L) shifts into K.
$ws.Columns("K").Delete()

# Remove the last data row (old row 6, "Normal hearing").
$ws.Rows("6").Delete()

# Row 2: Sex
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Sex"
$ws.Range("C2").Value = 99
$ws.Range("D2").Value = 100
$ws.Range("E2").Value = 51
$ws.Range("F2").Value = 100
$ws.Range("G2").Value = 48
$ws.Range("H2").Value = 100
$ws.Range("I2").Value = 4.371473494780586
$ws.Range("J2").Value = 0.0436
$ws.Range("K2").Value = "Significant"

# Row 3: Temporal_Bone_CT_scan_Findings_Pathological_Ear
$ws.Range("A3").Value = 16
$ws.Range("B3").Value = "Temporal_Bone_CT_scan_Findings_Pathological_Ear"
$ws.Range("C3").Value = 99
$ws.Range("D3").Value = 100
$ws.Range("E3").Value = 51
$ws.Range("F3").Value = 100
$ws.Range("G3").Value = 48
$ws.Range("H3").Value = 100
$ws.Range("I3").Value = 74.35668103448276
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = "Significant"

# Row 4: Surgical_Treatment
$ws.Range("A4").Value = 20
$ws.Range("B4").Value = "Surgical_Treatment"
$ws.Range("C4").Value = 51
$ws.Range("D4").Value = 100
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 100
$ws.Range("G4").Value = 48
$ws.Range("H4").Value = 100
$ws.Range("I4").Value = 33.30612244897959
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = "Significant"

# Row 5: One_Week_Post_surgical_Tonal_Audiometry_Pathological_Ear
$ws.Range("A5").Value = 27
$ws.Range("B5").Value = "One_Week_Post_surgical_Tonal_Audiometry_Pathological_Ear"
$ws.Range("C5").Value = 33
$ws.Range("D5").Value = 100
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 100
$ws.Range("G5").Value = 32
$ws.Range("H5").Value = 100
$ws.Range("I5").Value = 33
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = "Significant"
